# Update imputed values in the RandomForest result sheet (re-run of the
# imputation algorithm produced slightly different numeric results).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.8778
$ws.Range("D3").Value = -7.80979999999999
$ws.Range("B4").Value = 4.876800000000002
$ws.Range("C4").Value = -14.45680000000001
$ws.Range("D4").Value = -7.854999999999998
$ws.Range("C5").Value = -14.90040000000002
$ws.Range("B6").Value = 9.6822
$ws.Range("B7").Value = 6.111099999999998
$ws.Range("C8").Value = -11.71359999999999
$ws.Range("D9").Value = -7.929400000000003
$ws.Range("D11").Value = -8.609499999999995
$ws.Range("D14").Value = -7.570500000000001
$ws.Range("B16").Value = 9.367600000000007
$ws.Range("C16").Value = -12.0976
$ws.Range("D18").Value = -8.409899999999995
$ws.Range("B20").Value = 4.821700000000002
$ws.Range("E20").Value = 12.30739999999999
$ws.Range("C22").Value = -11.23359999999999
$ws.Range("D25").Value = -8.418099999999999
